$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row7 = New-Object "object[,]" 1,73
$row7[0,0] = "2024-09-12_E_e.dat"
$row7[0,1] = "transfer"
$row7[0,2] = 200
$row7[0,3] = "Blackman"
$row7[0,4] = 0.447
$row7[0,5] = 0.0028
$row7[0,6] = 0.0148
$row7[0,7] = 0.00055
$row7[0,8] = 10818242.33717204
$row7[0,9] = 377
$row7[0,10] = 13
$row7[0,11] = 0.01
$row7[0,12] = 1.353251525585472
$row7[0,13] = 0.01823642281480822
$row7[0,14] = 3.384760825824527
$row7[0,15] = 0.04561305003519604
$row7[0,16] = 268.9189189189189
$row7[0,17] = 42.68013515262911
$row7[0,18] = 0.597685403414816
$row7[0,19] = 0.3098666670904903
$row7[0,20] = 0.01836981553431876
$row7[0,21] = 0.3473690933134114
$row7[0,22] = 2.176582100597987
$row7[0,23] = 0.3282364826248091
$row7[0,24] = 2.523951193911398
$row7[0,25] = 7.689429199728545
$row7[0,26] = 0.2348580854810282
$row7[0,27] = 0.01013078014731138
$row7[0,28] = 0.3373026864360487
$row7[0,29] = 0.344974735785374
$row7[0,30] = 0.3298719555730663
$row7[0,31] = 0.3374450992921992
$row7[0,32] = 0.008089328022197521
$row7[0,33] = 2.650770589112001
$row7[0,34] = 2.70922255528005
$row7[0,35] = 2.581354759869761
$row7[0,36] = 2.646708965697791
$row7[0,37] = 0.06443938461543329
$row7[0,38] = 7.862530888593243
$row7[0,39] = 8.071431991609899
$row7[0,40] = 7.621923052448456
$row7[0,41] = 7.846747444082977
$row7[0,42] = 0.232870309309199
$row7[0,43] = 0.01943038952938295
$row7[0,44] = 0.01991282363412602
$row7[0,45] = 0.01886497187300008
$row7[0,46] = 0.01939564278916637
$row7[0,47] = 0.000526360150116652
$row7[0,48] = 2.302253037201362
$row7[0,49] = 2.359176934236709
$row7[0,50] = 2.235269316893389
$row7[0,51] = 2.297922443227667
$row7[0,52] = 0.06207740544018735
$row7[0,53] = 0.3486887782037926
$row7[0,54] = 0.3563837980805823
$row7[0,55] = 0.3405055953707706
$row7[0,56] = 0.3487865224701245
$row7[0,57] = 0.008845137064875557
$row7[0,58] = 1.454080806710048
$row7[0,59] = 1.490033447924836
$row7[0,60] = 1.411775459226118
$row7[0,61] = 1.451328557431822
$row7[0,62] = 0.03920150396585813
$row7[0,63] = 2.15540550673144
$row7[0,64] = 2.219596836482521
$row7[0,65] = 2.082332307800371
$row7[0,66] = 2.151459300582226
$row7[0,67] = 0.07072926319807023
$row7[0,68] = 3.630071685653816
$row7[0,69] = 0.09805103665383658
$row7[0,70] = 3.636955628042732
$row7[0,71] = 3.726880589713171
$row7[0,72] = 3.531141239380172
$ws.Range("A7:BU7").Value = $row7

$row8 = New-Object "object[,]" 1,73
$row8[0,0] = "2024-09-18_F_e.dat"
$row8[0,1] = "transfer"
$row8[0,2] = 200
$row8[0,3] = "Blackman"
$row8[0,4] = 0.6889999999999999
$row8[0,5] = 0.009900000000000001
$row8[0,6] = 0.0186
$row8[0,7] = 0.00014
$row8[0,8] = 12127808.20404177
$row8[0,9] = 377
$row8[0,10] = 13
$row8[0,11] = 0.01
$row8[0,12] = 0.7663633636591822
$row8[0,13] = 0.01823642281480822
$row8[0,14] = 1.709852237829378
$row8[0,15] = 0.0406877335719948
$row8[0,16] = 213.9784946236559
$row8[0,17] = 136.2503920113797
$row8[0,18] = 6.091112189062183
$row8[0,19] = 0.8170872608259792
$row8[0,20] = 0.003384662392636809
$row8[0,21] = 0.03696083600119154
$row8[0,22] = 0.7362064422805231
$row8[0,23] = 0.820471923218616
$row8[0,24] = 0.7731672782817146
$row8[0,25] = 0.9423445902312773
$row8[0,26] = 0.4698072715794814
$row8[0,27] = 0.01653871829921932
$row8[0,28] = 0.351763617729176
$row8[0,29] = 0.3669817794881653
$row8[0,30] = 0.3376683775226238
$row8[0,31] = 0.3521845749111452
$row8[0,32] = 0.01373476543499438
$row8[0,33] = 0.9557550836038384
$row8[0,34] = 0.9891662673434409
$row8[0,35] = 0.9236846561540357
$row8[0,36] = 0.9549304444654596
$row8[0,37] = 0.03260395484408993
$row8[0,38] = 2.711787530709346
$row8[0,39] = 2.873109813708927
$row8[0,40] = 2.572900857329363
$row8[0,41] = 2.715749851987404
$row8[0,42] = 0.1440443472497804
$row8[0,43] = 0.003631281234409604
$row8[0,44] = 0.003779725136137824
$row8[0,45] = 0.003485725796355147
$row8[0,46] = 0.003630970176366508
$row8[0,47] = 0.0001448439314680616
$row8[0,48] = 0.7898508346236919
$row8[0,49] = 0.8221388795426557
$row8[0,50] = 0.7581913840172039
$row8[0,51] = 0.7897833767657624
$row8[0,52] = 0.0315054717437733
$row8[0,53] = 0.1653459088236783
$row8[0,54] = 0.1715149702597857
$row8[0,55] = 0.1578411112170259
$row8[0,56] = 0.1651470676996971
$row8[0,57] = 0.006795520846426978
$row8[0,58] = 0.6490102675788103
$row8[0,59] = 0.6755408578834345
$row8[0,60] = 0.6229962442106776
$row8[0,61] = 0.6489548759671027
$row8[0,62] = 0.02588764734473556
$row8[0,63] = 0.9231716820827645
$row8[0,64] = 0.9786291849319789
$row8[0,65] = 0.8705223432078169
$row8[0,66] = 0.9227996727485688
$row8[0,67] = 0.05271334586028383
$row8[0,68] = 1.447899259725191
$row8[0,69] = 0.05775856968577506
$row8[0,70] = 1.448022845318822
$row8[0,71] = 1.507215900929803
$row8[0,72] = 1.38998231496444
$ws.Range("A8:BU8").Value = $row8
